# -----------------------------------------------------------------------
# Rename existing sheet to "gen", add "lines" and "nodes" after it so the
# final tab order is gen, nodes, lines.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$gen = $wb.Worksheets.Item(1)
$gen.Name = "gen"

$lines = $wb.Worksheets.Add($null, $gen)
$lines.Name = "lines"

$nodes = $wb.Worksheets.Add($null, $gen)
$nodes.Name = "nodes"

# ---------------------------------------------------------------------------
# Sheet "gen" : NODE | PMIN | PMAX | CV
# ---------------------------------------------------------------------------
$gen.Range("B1").Value = "PMIN"
$gen.Range("C1").Value = "PMAX"
$gen.Range("A1").Value = "NODE"
$gen.Range("D1").Value = "CV"

# ---------------------------------------------------------------------------
# Sheet "nodes"
# ---------------------------------------------------------------------------
$nodes.Range("G2").Value = "NODE I"
$nodes.Range("A2").Value = "NUM NODE"
$nodes.Range("B1").Value = "GEN"
$nodes.Range("D1").Value = "LOAD"
$nodes.Range("D2").Value = "P"
$nodes.Range("H2").Value = "NODE J"
$nodes.Range("I2").Value = "X"
$nodes.Range("J2").Value = "FMAX"
$nodes.Range("G1").Value = "LINES"

$nodes.Range("B2").Value = "PMIN"
$nodes.Range("C2").Value = "PMAX"

$nodes.Range("A3").Value = 1
$nodes.Range("B3").Value = 0
$nodes.Range("C3").Value = 100
$nodes.Range("D3").Value = 50

$nodes.Range("A4").Value = 2
$nodes.Range("B4").Value = 0
$nodes.Range("C4").Value = 100
$nodes.Range("D4").Value = 20

$nodes.Range("A5").Value = 3
$nodes.Range("B5").Value = 0
$nodes.Range("C5").Value = 100
$nodes.Range("D5").Value = 10

$nodes.Range("G3").Value = 1
$nodes.Range("H3").Value = 2
$nodes.Range("J3").Value = 200

$nodes.Range("G4").Value = 2
$nodes.Range("H4").Value = 3
$nodes.Range("J4").Value = 200

$nodes.Range("G5").Value = 3
$nodes.Range("H5").Value = 1
$nodes.Range("J5").Value = 200

# "0.4" stored as text (not a number) in I3:I5
$nodes.Range("I3:I5").NumberFormat = "@"
$nodes.Range("I3").Value = "0.4"
$nodes.Range("I4").Value = "0.4"
$nodes.Range("I5").Value = "0.4"
$nodes.Range("I3:I5").NumberFormat = "General"

# Merge header cells
$nodes.Range("B1:C1").Merge()
$nodes.Range("G1:J1").Merge()

# Alignment: center the whole data block
$nodes.Range("A1:D5").HorizontalAlignment = -4108
$nodes.Range("G1:J5").HorizontalAlignment = -4108
$nodes.Range("G3:H5").HorizontalAlignment = -4131
$nodes.Range("J3:J5").HorizontalAlignment = -4131

# Bold + colored fill for the "GEN" / "LOAD" / "LINES" banner cells
$nodes.Range("B1:C1").Font.Bold = $true
$nodes.Range("B1:C1").Interior.ThemeColor = 6
$nodes.Range("B1:C1").Interior.TintAndShade = 0.39997558519241921

$nodes.Range("D1").Font.Bold = $true
$nodes.Range("D1").Interior.ThemeColor = 5
$nodes.Range("D1").Interior.TintAndShade = 0.39997558519241921

$nodes.Range("G1:J1").Interior.Color = 49407

# ---------------------------------------------------------------------------
# Defined names (workbook scoped)
# ---------------------------------------------------------------------------
$wb.Names.Add("LOAD", "=nodes!`$D`$1:`$D`$5")
$wb.Names.Add("PMAX", "=nodes!`$C`$2:`$C`$5")
$wb.Names.Add("PMIN", "=nodes!`$B`$2:`$B`$5")

# ---------------------------------------------------------------------------
# Final selection / active sheet state
# ---------------------------------------------------------------------------
$gen.Range("A2").Select()
$nodes.Select()
$nodes.Range("K13:L13").Select()
